$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix punctuation in the description of B7: a period before "hacen" should
# be a comma, since the sentence continues rather than starting a new one.
$cell = $ws.Range("B7")
$text = $cell.Value2
$fixed = $text.Replace('(con z). hacen', '(con z), hacen')
$cell.Value = $fixed

# Leave the selection on B7 alone (single cell) instead of the whole row
# A7:XFD7 that was previously selected.
$ws.Range("B7").Select()
